$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets a distinct value
$ws.Range("C2").Value = 7318

# Rows 3 through 252 all get the same updated value
$ws.Range("C3:C252").Value = 7310
